$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update D2, E2
$ws.Range("D2").Value = '33.875.88'
$ws.Range("E2").Value = '  -2.32%  '

# Row 3: update D3, E3
$ws.Range("D3").Value = '1.768.00'
$ws.Range("E3").Value = '  -1.47%  '

# Row 4: update E4
$ws.Range("E4").Value = '  -0.08%  '

# Row 5: update E5
$ws.Range("E5").Value = '  -2.32%  '

# Row 6: update D6, E6
$ws.Range("D6").Value = '0.546'
$ws.Range("E6").Value = '  -1.51%  '

# Row 7: update E7
$ws.Range("E7").Value = '  -0.07%  '

# Row 8: update D8, E8
$ws.Range("D8").Value = '30.83'
$ws.Range("E8").Value = '  -5.99%  '

# Row 9: update E9
$ws.Range("E9").Value = '  -0.33%  '

# Row 10: update D10, E10
$ws.Range("D10").Value = '0.0701'
$ws.Range("E10").Value = '  +4.82%  '

# Row 11: update D11, E11
$ws.Range("D11").Value = '0.0921'
$ws.Range("E11").Value = '  -1.64%  '

# Row 12: update D12
$ws.Range("D12").Value = '2.024.63'

# Row 13: update D13, E13
$ws.Range("D13").Value = '1.763.96'
$ws.Range("E13").Value = '  -1.54%  '

# Row 14: update D14, E14
$ws.Range("D14").Value = '10.44'
$ws.Range("E14").Value = '  -6.00%  '

# Row 15: update E15
$ws.Range("E15").Value = '  -2.29%  '

# Row 16: update D16, E16
$ws.Range("D16").Value = '33.905.68'
$ws.Range("E16").Value = '  -2.28%  '

# Row 17: update E17
$ws.Range("E17").Value = '  -2.60%  '

# Row 18: update D18, E18
$ws.Range("D18").Value = '67.48'
$ws.Range("E18").Value = '  -2.49%  '

# Row 19: update D19, E19
$ws.Range("D19").Value = '242.61'
$ws.Range("E19").Value = '  -4.76%  '

# Row 20: update D20, E20
$ws.Range("D20").Value = '0.0₃0771'
$ws.Range("E20").Value = '  +0.87%  '

# Row 21: update D21, E21
$ws.Range("D21").Value = '0.999'
$ws.Range("E21").Value = '  -0.21%  '

# Row 22: update D22, E22
$ws.Range("D22").Value = '10.48'
$ws.Range("E22").Value = '  +0.69%  '

# Row 23: update D23, E23
$ws.Range("D23").Value = '4.02'
$ws.Range("E23").Value = '  -5.15%  '

# Row 24: update E24
$ws.Range("E24").Value = '  -1.41%  '

# Row 25: update D25, E25
$ws.Range("D25").Value = '157.19'
$ws.Range("E25").Value = '  -0.88%  '

# Row 26: update D26, E26
$ws.Range("D26").Value = '16.28'
$ws.Range("E26").Value = '  -0.88%  '

# Row 27: update D27, E27
$ws.Range("D27").Value = '6.94'
$ws.Range("E27").Value = '  -2.16%  '

# Row 28: update E28
$ws.Range("E28").Value = '  -2.74%  '

# Row 29: update E29
$ws.Range("E29").Value = '  +0.06%  '

# Row 30: update D30, E30
$ws.Range("D30").Value = "'" + '0.0520'
$ws.Range("E30").Value = '  +0.25%  '

# Row 31: update D31, E31
$ws.Range("D31").Value = '3.68'
$ws.Range("E31").Value = '  -2.33%  '

# Row 32: update E32
$ws.Range("E32").Value = '  +0.12%  '

# Row 33: update E33
$ws.Range("E33").Value = '  -2.87%  '

# Row 34: update E34
$ws.Range("E34").Value = '  -3.76%  '

# Row 35: update D35, E35
$ws.Range("D35").Value = '1.394.62'
$ws.Range("E35").Value = '  -3.88%  '

# Row 36: update D36, E36
$ws.Range("D36").Value = '1.05'
$ws.Range("E36").Value = '  -1.57%  '

# Row 37: update D37, E37
$ws.Range("D37").Value = '0.631'
$ws.Range("E37").Value = '  +0.29%  '

# Row 38: update E38
$ws.Range("E38").Value = '  -2.30%  '

# Row 39: update D39, E39
$ws.Range("D39").Value = '0.923'
$ws.Range("E39").Value = '  +2.53%  '

# Row 40: update E40
$ws.Range("E40").Value = '  -0.50%  '

# Row 41: update D41, E41
$ws.Range("D41").Value = '78.37'
$ws.Range("E41").Value = '  -5.80%  '

# Row 42: update E42
$ws.Range("E42").Value = '  -5.47%  '

# Row 43: update D43, E43
$ws.Range("D43").Value = '2.09'
$ws.Range("E43").Value = '  +1.04%  '

# Row 44: update E44
$ws.Range("E44").Value = '  -1.32%  '

# Row 45: update B45, C45, D45, E45
$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").Value = '1.04'
$ws.Range("E45").Value = '  -1.63%  '

# Row 46: update B46, C46, D46, E46
$ws.Range("B46").Value = 'Kaspa'
$ws.Range("C46").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D46").Value = '0.0487'
$ws.Range("E46").Value = '  -4.07%  '

# Row 47: update D47, E47
$ws.Range("D47").Value = '1.920.08'
$ws.Range("E47").Value = '  -1.86%  '

# Row 48: update D48, E48
$ws.Range("D48").Value = '103.73'
$ws.Range("E48").Value = '  -0.40%  '

# Row 49: update E49
$ws.Range("E49").Value = '  -0.68%  '

# Row 50: update D50, E50
$ws.Range("D50").Value = '11.77'
$ws.Range("E50").Value = '  -1.26%  '

# Row 51: update E51
$ws.Range("E51").Value = '  -2.73%  '
